# Updating with July charges
$wb = $excel.ActiveWorkbook
$wsCodes = $wb.Worksheets.Item("CL Codes")
$wsAlt = $wb.Worksheets.Item("eto_use_alt")

# NOTE: cell writes below are deliberately ordered to match the order in
# which brand-new shared strings are first introduced in the saved file
# (CL015 -> "CL015, CL006, CL008" -> "Guillaume Hoareau/Austin Johnson" ->
# "Hoareau/Johnson"), so the regenerated sharedStrings table lines up.

# --- "CL Codes" sheet: new code row (row 18), column A first ---
$wsCodes.Range("A18").Value = "CL015"

# --- "eto_use_alt" sheet: new July usage rows (31-33) ---
$wsAlt.Range("A31").Value = [DateTime]"2024-07-08"
$wsAlt.Range("B31").Value = "CL015, CL006, CL008"

# --- back to "CL Codes" for the remaining new-string columns ---
$wsCodes.Range("B18").Value = "Guillaume Hoareau/Austin Johnson"
$wsCodes.Range("C18").Value = "Hoareau/Johnson"
$wsCodes.Range("D18").Value = "emergency medicine"

# --- remaining "eto_use_alt" rows (reuse already-existing strings) ---
$wsAlt.Range("A32").Value = [DateTime]"2024-07-09"
$wsAlt.Range("B32").Value = "CL015"

$wsAlt.Range("A33").Value = [DateTime]"2024-07-17"
$wsAlt.Range("B33").Value = "CL008"

# --- Selections / active sheet to match saved view state ---
$wsAlt.Range("B35").Select()
$wsCodes.Range("D22").Select()
$wsCodes.Activate()
